$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(2, 4, 5, 6, 7, 9, 10, 11, 12, 13)

$data = @(
    @(2, 0.6374623883579886, 0.01519605238401311, 0.03918841062479927, 4.864748571270326, 0.002662847161076502, 2.748841524720092, 0.1121708070941665, 1.346030194322822, 0.5487130501106776, 0.2917769988649752),
    @(3, 0.6332289875409884, 0.01328457339317879, 0.03809945734230524, 4.817285731785901, 0.002667864416992662, 2.72863149191835, 0.1109687631005603, 1.303687818532097, 0.5432273343945297, 0.2893252131999198),
    @(4, 0.6310672385716174, 0.01210861245044015, 0.03742010295826503, 4.789513360369313, 0.002671108453819828, 2.716733177454444, 0.1102052986379416, 1.278684324657945, 0.5401563956606594, 0.2879977910819136),
    @(5, 0.6302965998241916, 0.01162876167101956, 0.0371405238830409, 4.77853990042928, 0.002672471654427505, 2.712012688740103, 0.1098877446877147, 1.268744987606908, 0.5389797116520185, 0.2875016629578475),
    @(6, 0.6301753048851566, 0.01154904265462164, 0.03709393360128388, 4.776738528779617, 0.002672700506968562, 2.711236588974757, 0.1098346249829856, 1.267109644035202, 0.5387888384792774, 0.2874219890578011),
    @(7, 0.6310563985569644, 0.01210214367102225, 0.03741634358690504, 4.789363976255373, 0.002671126671286116, 2.716668996686138, 0.1102010421154223, 1.278549268396773, 0.5401402238521342, 0.28799091863333),
    @(8, 0.6359120436452201, 0.01453741770528438, 0.03881514405236963, 4.848098721022211, 0.002664543274680953, 2.741766941514257, 0.1117615860920154, 1.331223799131152, 0.5467598559225451, 0.2908947113177405),
    @(9, 0.6488962564131953, 0.0192974887023567, 0.0414748278791599, 4.974179287196563, 0.002652923730333173, 2.795051508212026, 0.1146224484614855, 1.442438042531222, 0.5621021559393284, 0.2979998171700267),
    @(10, 0.6605360475183772, 0.02278942924792204, 0.0433807036956626, 5.073513976841355, 0.002645164874775772, 2.83670567542012, 0.1166059384460638, 1.52902250636194, 0.5748189064502611, 0.3040792565784827),
    @(11, 0.6662853715909307, 0.02437777238441896, 0.04423784696952993, 5.12017274079642, 0.002641802259232421, 2.856205390323382, 0.1174832568765094, 1.569481828201077, 0.5809191429128333, 0.3070314997247863),
    @(12, 0.6685276177808817, 0.02497928185291443, 0.04456105247916575, 5.138053529298702, 0.002640552787025489, 2.863669038000907, 0.1178119355320799, 1.584957512729176, 0.5832745586738071, 0.308176257589345),
    @(13, 0.6680418181260563, 0.02484973334098584, 0.04449150494540888, 5.134193137303157, 0.002640820823260779, 2.862058066190883, 0.117741305479349, 1.581617662049297, 0.5827652583114826, 0.3079285217895134),
    @(14, 0.6664685390935574, 0.02442725780878163, 0.04426446462758804, 5.121639549147204, 0.002641698986671686, 2.856817833554814, 0.1175103681859238, 1.570751920275882, 0.5811120144563233, 0.3071251425830965),
    @(15, 0.6655133317177331, 0.02416848616788769, 0.04412521787590507, 5.113977757458031, 0.002642239992857361, 2.853618405456018, 0.1173684525215801, 1.564116491187576, 0.5801052672078697, 0.3066365396995536),
    @(16, 0.6601694381365633, 0.02268562932437135, 0.04332449260981264, 5.070494367503045, 0.002645387977953224, 2.835442436156725, 0.1165481051621082, 1.526400002271686, 0.5744265905668868, 0.3038900737937666),
    @(17, 0.6570073313927765, 0.02177594755746526, 0.0428307795196794, 5.04419588482736, 0.00264736182846335, 2.824433405767678, 0.1160384854787413, 1.503536989284527, 0.5710236994184186, 0.3022529916289187),
    @(18, 0.6552313478782708, 0.02125270410304836, 0.04254587932850207, 5.029208105575691, 0.002648512854520835, 2.818153150228113, 0.1157430114024613, 1.490487622253568, 0.5690961231757399, 0.3013289581203757),
    @(19, 0.654637384362502, 0.0210755378157117, 0.04244925615902595, 5.024157255022146, 0.002648905275596271, 2.8160356596009, 0.1156425630991382, 1.486086633425117, 0.5684485746573955, 0.3010191155402353),
    @(20, 0.6573395175157941, 0.02187278621280342, 0.04288343202676081, 5.046981073702028, 0.002647150082510005, 2.825599967737702, 0.1160929786667886, 1.505960356476237, 0.5713828712469819, 0.3024254433683708),
    @(21, 0.6669288842413721, 0.02455134758280053, 0.0443311888953648, 5.125321080297624, 0.002641440401861792, 2.858354856177812, 0.1175782958221951, 1.573939251519732, 0.5815963800375243, 0.307360387263202),
    @(22, 0.673575450241799, 0.02630219472391104, 0.04526937997463243, 5.177757533864792, 0.00263784791221321, 2.880225761083807, 0.1185284094697643, 1.619268871505199, 0.5885360816814824, 0.3107418972582821),
    @(23, 0.6699934160507439, 0.02536769125444494, 0.04476936850974411, 5.149657846214239, 0.002639752603998964, 2.868510322706058, 0.1180231868963002, 1.594992938124818, 0.5848080069359582, 0.3089228373611661),
    @(24, 0.6571892054721218, 0.02182900625656004, 0.04285963112738145, 5.045721480021768, 0.002647245762198058, 2.825072412508888, 0.1160683500477671, 1.504864456660613, 0.5712203999855774, 0.3023474245241289),
    @(25, 0.6450141124998652, 0.01801099314594978, 0.04076396698524754, 4.93889814968739, 0.002655929871056249, 2.780198446950322, 0.1173684525215801, 1.564116491187576, 0.5801052672078697, 0.3066365396995536)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $colIndexes.Count; $i++) {
        $ws.Cells.Item($r, $colIndexes[$i]).Value = $row[$i + 1]
    }
}

"done"